$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.195.99"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "3.600.02"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'603.14"
$ws.Range("D6").Value = "'196.05"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").Value = "'0.647"
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").Value = "'53.76"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "'9.55"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "4.180.71"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Value = "'13.11"
$ws.Range("E15").Value = "  +3.52%  "
$ws.Range("D16").Value = "'597.21"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "70.437.70"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.615.69"
$ws.Range("E18").Value = "  +2.46%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'19.06"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "'0.995"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").Value = "'17.83"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").Value = "'5.17"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("D24").Value = "'102.11"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("D25").Value = "'4.61"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("D27").Value = "'10.74"
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("D28").Value = "'9.62"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("D29").Value = "'33.78"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").Value = "'4.78"
$ws.Range("E30").Value = "  +6.19%  "
$ws.Range("D31").Value = "'7.15"
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").Value = "'12.28"
$ws.Range("E32").Value = "  -3.62%  "
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("D34").Value = "'63.26"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").Value = "0.0₃0899"
$ws.Range("E35").Value = "  +8.37%  "
$ws.Range("D36").Value = "3.905.19"
$ws.Range("E36").Value = "  +4.37%  "
$ws.Range("D37").Value = "'3.09"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "'521.94"
$ws.Range("E39").Value = "  +5.41%  "
$ws.Range("D40").Value = "'36.90"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").Value = "'0.389"
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("D42").Value = "'3.52"
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").Value = "'3.45"
$ws.Range("E45").Value = "  +3.71%  "
$ws.Range("D46").Value = "'2.85"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "'8.62"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "'0.000251"
$ws.Range("E50").Value = "  +2.27%  "
$ws.Range("E51").Value = "  +0.05%  "
